# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13701
$ws1.Range("F3").Value = 325
$ws1.Range("F4").Value = 663
$ws1.Range("F5").Value = 230
$ws1.Range("F6").Value = 481
$ws1.Range("F7").Value = 1399
$ws1.Range("F8").Value = 131

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13701
$ws4.Range("F3").Value = 325
$ws4.Range("F4").Value = 663
$ws4.Range("F5").Value = 230
$ws4.Range("F8").Value = 481
$ws4.Range("F9").Value = 1399
$ws4.Range("F11").Value = 131
